# fibonacciUpTo50.xlsx -> fibonacciUpTo100.xlsx
# "changed int to unsigned long long, implemented tests for fibonacci and fibonacciSum"
# Extends the Fibonacci table from 50 terms (rows 2-52) to 100 terms (rows 2-102),
# widens columns B/C, moves the view down to the new data, and highlights the
# last 10 rows (91-100, i.e. sheet rows 93-102) in red.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for B (value) and C (running sum) -----------------------
$ws.Columns.Item(2).ColumnWidth = 18.5703125
$ws.Columns.Item(3).ColumnWidth = 27.42578125

# --- Extend the Fibonacci sequence from row 53 (index 51) to row 102 (index 100) --
for ($r = 53; $r -le 102; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Formula = "=`$B$($r-2) + `$B$($r-1)"
    $ws.Cells.Item($r, 3).Formula = "=SUM(`$B`$2:`$B$r)"
}

# --- Highlight the final ten rows (91..100 -> sheet rows 93..102) red ------
$ws.Range("A93:C93").Interior.Color = 65535
$ws.Range("A94:C102").Interior.Color = 255

# --- Refresh the used-range dimension / view to match the new data ---------
$ws.Range("F92").Select()
$excel.ActiveWindow.ScrollRow = 78
